# Auto-generated Excel COM-interop edit script
# Applies the commit diff: updates 'want-to-go' counts (F column) across
# sheets 1 (å±è§), 3 (æ¬å°çæ´»), 4 (å¨é¨ç±»å);
# and on sheet 2 (æ¼åº) removes the oldest event row (old row 2) and
# shifts every later event up by one row, then deletes the now-empty last row.

$wb = $excel.ActiveWorkbook
$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsPerform = $wb.Worksheets.Item(2)   # 演出
$wsLocal   = $wb.Worksheets.Item(3)   # 本地生活
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1 (展览): update "想去人数" (F column) counts ---
$wsExhibit.Range("F2").Value = 309
$wsExhibit.Range("F3").Value = 1255
$wsExhibit.Range("F4").Value = 361
$wsExhibit.Range("F5").Value = 328
$wsExhibit.Range("F6").Value = 3835
$wsExhibit.Range("F9").Value = 2213
$wsExhibit.Range("F10").Value = 328
$wsExhibit.Range("F12").Value = 735
$wsExhibit.Range("F13").Value = 158
$wsExhibit.Range("F14").Value = 160
$wsExhibit.Range("F15").Value = 2112
$wsExhibit.Range("F17").Value = 8
$wsExhibit.Range("F19").Value = 330

# --- Sheet 3 (本地生活): update "想去人数" (F column) counts ---
$wsLocal.Range("F2").Value = 6394
$wsLocal.Range("F3").Value = 813
$wsLocal.Range("F4").Value = 2076
$wsLocal.Range("F5").Value = 308

# --- Sheet 4 (全部类型): update "想去人数" (F column) counts ---
$wsAll.Range("F2").Value = 6394
$wsAll.Range("F3").Value = 813
$wsAll.Range("F4").Value = 2076
$wsAll.Range("F5").Value = 308
$wsAll.Range("F6").Value = 32
$wsAll.Range("F10").Value = 309
$wsAll.Range("F11").Value = 1255
$wsAll.Range("F12").Value = 361
$wsAll.Range("F16").Value = 328
$wsAll.Range("F17").Value = 3835
$wsAll.Range("F18").Value = 126
$wsAll.Range("F24").Value = 2213
$wsAll.Range("F25").Value = 328
$wsAll.Range("F28").Value = 735
$wsAll.Range("F29").Value = 158
$wsAll.Range("F30").Value = 160
$wsAll.Range("F32").Value = 2112
$wsAll.Range("F36").Value = 8
$wsAll.Range("F38").Value = 330

# --- Sheet 2 (演出): drop the 2024-10-19 event (old row 2) and shift all
# later events (old rows 3-23) up by one row; column A (row index) is left
# untouched since it already holds the correct sequential value per row. ---
$wsPerform.Range("B2").Value = '2024-10-20'
$wsPerform.Range("C2").Value = '广州·ROOKiEZ is PUNK`D 「Reignite Youth （重燃青春）」2024 CHINA Tour '
$wsPerform.Range("D2").Value = '南洲路154号侨建大厦2F SDlivehouse'
$wsPerform.Range("E2").Value = '2024.10.20 20:00-10.20 21:30'
$wsPerform.Range("F2").Value = 32
$wsPerform.Range("G2").Value = 259
$wsPerform.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=92075'
$wsPerform.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202409/kAxVF2Jw1725542237304.jpeg'
$wsPerform.Range("B3").Value = '2024-10-25'
$wsPerform.Range("C3").Value = '广州·东方PartyNight×东方同人only-游剧天P2'
$wsPerform.Range("D3").Value = '南洲路158号2F SD Livehouse'
$wsPerform.Range("E3").Value = '2024.10.25 19:00-10.25 22:00'
$wsPerform.Range("F3").Value = 39
$wsPerform.Range("G3").Value = 149
$wsPerform.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=93136'
$wsPerform.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202409/Oc5h2el91727671722939.jpeg'
$wsPerform.Range("B4").Value = '2024-10-25'
$wsPerform.Range("C4").Value = '广州·新生代流媒体小天后野田爱实 2024 巡演'
$wsPerform.Range("D4").Value = '南洲路158号2F SD Livehouse'
$wsPerform.Range("E4").Value = '2024.10.25 20:00-10.25 22:00'
$wsPerform.Range("F4").Value = 14
$wsPerform.Range("G4").Value = 280
$wsPerform.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=91823'
$wsPerform.Range("I4").Value = '//i0.hdslb.com/bfs/openplatform/202409/oN7FyQ8v1725347758464.jpeg'
$wsPerform.Range("B5").Value = '2024-10-27'
$wsPerform.Range("C5").Value = '广州·《次元共鸣》 ACG乐队番主题演唱会 后乐园乐队'
$wsPerform.Range("D5").Value = '南洲路158号2F SD Livehouse'
$wsPerform.Range("E5").Value = '2024.10.27 20:00-10.27 22:00'
$wsPerform.Range("F5").Value = 23
$wsPerform.Range("G5").Value = 160
$wsPerform.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=93166'
$wsPerform.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202410/ADwxxu0H1728383860332.jpeg'
$wsPerform.Range("B6").Value = '2024-10-27'
$wsPerform.Range("C6").Value = '广州·卡农·世界经典音乐之旅交响音乐会'
$wsPerform.Range("D6").Value = '东风中路299号 广州中山纪念堂'
$wsPerform.Range("E6").Value = '2024.10.27 19:30-10.27 21:00'
$wsPerform.Range("F6").Value = 34
$wsPerform.Range("G6").Value = 75
$wsPerform.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=91040'
$wsPerform.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202408/WEqD8aj31724134831558.jpeg'
$wsPerform.Range("B7").Value = '2024-10-27'
$wsPerform.Range("C7").Value = '广州·混合理论——致敬林肯公园世界巡回演唱会'
$wsPerform.Range("D7").Value = '龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House'
$wsPerform.Range("E7").Value = '2024.10.27 20:00-10.27 21:30'
$wsPerform.Range("F7").Value = 9
$wsPerform.Range("G7").Value = 280
$wsPerform.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=92921'
$wsPerform.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202409/ZRv3pMgC1727244329886.jpeg'
$wsPerform.Range("B8").Value = '2024-11-03'
$wsPerform.Range("C8").Value = '广州·majiko巡演-2024'
$wsPerform.Range("D8").Value = '龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House'
$wsPerform.Range("E8").Value = '2024.11.03 20:00-11.03 21:40'
$wsPerform.Range("F8").Value = 125
$wsPerform.Range("G8").Value = 480
$wsPerform.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=92291'
$wsPerform.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202409/5t950dFT1726035772853.jpeg'
$wsPerform.Range("B9").Value = '2024-11-08'
$wsPerform.Range("C9").Value = '广州·HAG·CHINA TOUR 2024 - Phase 2 -1st Lve in Guangzhou「 声 」'
$wsPerform.Range("D9").Value = '机场路1733号 久米空间LIVEHOUSE'
$wsPerform.Range("E9").Value = '2024.11.08 19:30-11.10 21:00'
$wsPerform.Range("F9").Value = 91
$wsPerform.Range("G9").Value = 480
$wsPerform.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=93120'
$wsPerform.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202410/3q1wqY2V1728354279220.jpeg'
$wsPerform.Range("B10").Value = '2024-11-08'
$wsPerform.Range("C10").Value = '广州·「心随歌行」KOKIA 2024 中国巡演'
$wsPerform.Range("D10").Value = '广州大道中1229号 广东艺术剧院'
$wsPerform.Range("E10").Value = '2024.11.08 19:30-11.08 21:30'
$wsPerform.Range("F10").Value = 45
$wsPerform.Range("G10").Value = '不可售'
$wsPerform.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=90392'
$wsPerform.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202408/FDsbokRk1722914443578.jpeg'
$wsPerform.Range("B11").Value = '2024-11-09'
$wsPerform.Range("C11").Value = '广州·「多厨狂喜」白金交响乐团二次元交响音乐会'
$wsPerform.Range("D11").Value = '广州大道中1229号 广东艺术剧院'
$wsPerform.Range("E11").Value = '2024.11.09 19:30-11.09 21:30'
$wsPerform.Range("F11").Value = 83
$wsPerform.Range("G11").Value = 188
$wsPerform.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=93183'
$wsPerform.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202410/FixLtKw71728419735084.jpeg'
$wsPerform.Range("B12").Value = '2024-11-10'
$wsPerform.Range("C12").Value = '广州·平田雄也&小池亮介2024粉丝见面会'
$wsPerform.Range("D12").Value = '金花街道中山七路333号1906科技圆区3号楼109-1铺、110-1铺、111-1铺 音乐唐人馆'
$wsPerform.Range("E12").Value = '2024.11.10 13:00-11.10 18:00'
$wsPerform.Range("F12").Value = 224
$wsPerform.Range("G12").Value = 480
$wsPerform.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=92655'
$wsPerform.Range("I12").Value = '//i2.hdslb.com/bfs/openplatform/202409/UkhOeOwe1726658317935.jpeg'
$wsPerform.Range("B13").Value = '2024-11-17'
$wsPerform.Range("C13").Value = '广州·“法国姐姐”乔伊丝·乔纳森《小意思》2024巡回演唱会'
$wsPerform.Range("D13").Value = '东风中路299号 广州中山纪念堂'
$wsPerform.Range("E13").Value = '2024.11.17 19:30-11.17 21:00'
$wsPerform.Range("F13").Value = 7
$wsPerform.Range("G13").Value = 280
$wsPerform.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=91814'
$wsPerform.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202408/bnKPQEEd1725008600562.jpeg'
$wsPerform.Range("B14").Value = '2024-11-24'
$wsPerform.Range("C14").Value = '广州·【限时早鸟8折】奇妙人声之旅 · RESOUND理想人声阿卡贝拉音乐会 '
$wsPerform.Range("D14").Value = '广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）'
$wsPerform.Range("E14").Value = '2024.11.24 20:00-11.24 21:30'
$wsPerform.Range("F14").Value = 2
$wsPerform.Range("G14").Value = 144
$wsPerform.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=90940'
$wsPerform.Range("I14").Value = '//i2.hdslb.com/bfs/openplatform/202408/q7p66BEy1724037045076.jpeg'
$wsPerform.Range("B15").Value = '2024-11-27'
$wsPerform.Range("C15").Value = '广州·三重唱Ohashi Trio（大桥トリ才） 2024年巡演'
$wsPerform.Range("D15").Value = '人民北路875号（广州市少年宫内） 广州蓓蕾剧院'
$wsPerform.Range("E15").Value = '2024.11.27 19:30-11.27 21:00'
$wsPerform.Range("F15").Value = 7
$wsPerform.Range("G15").Value = 380
$wsPerform.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=91847'
$wsPerform.Range("I15").Value = '//i2.hdslb.com/bfs/openplatform/202409/ggAAQH8D1725369168304.jpeg'
$wsPerform.Range("B16").Value = '2024-12-08'
$wsPerform.Range("C16").Value = '广州·梁祝之父：何占豪指挥《梁祝》65周年大型东方交响音乐会'
$wsPerform.Range("D16").Value = '东风中路299号 广州中山纪念堂'
$wsPerform.Range("E16").Value = '2024.12.08 19:30-12.08 21:10'
$wsPerform.Range("F16").Value = 7
$wsPerform.Range("G16").Value = 70
$wsPerform.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=92833'
$wsPerform.Range("I16").Value = '//i1.hdslb.com/bfs/openplatform/202409/y8ck801y1726297263642.jpeg'
$wsPerform.Range("B17").Value = '2024-12-11'
$wsPerform.Range("C17").Value = '广州·安田丽（安田レイ）「无形之线」2024巡演'
$wsPerform.Range("D17").Value = '新滘中路88号海珠同创汇东一街11号 声音共和Livehouse'
$wsPerform.Range("E17").Value = '2024.12.11 20:00-12.11 22:00'
$wsPerform.Range("F17").Value = 40
$wsPerform.Range("G17").Value = 320
$wsPerform.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=91909'
$wsPerform.Range("I17").Value = '//i0.hdslb.com/bfs/openplatform/202409/2821JdMa1725357077006.jpeg'
$wsPerform.Range("B18").Value = '2024-12-14'
$wsPerform.Range("C18").Value = '广州·变形金刚音乐会40周年特变版'
$wsPerform.Range("D18").Value = '广州大道中1229号 广东艺术剧院'
$wsPerform.Range("E18").Value = '2024.12.14 19:30-12.14 21:30'
$wsPerform.Range("F18").Value = 51
$wsPerform.Range("G18").Value = 171
$wsPerform.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=90033'
$wsPerform.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202407/RAV6qAVB1722168641097.jpeg'
$wsPerform.Range("B19").Value = '2024-12-20'
$wsPerform.Range("C19").Value = '广州·小野丽莎2024“倾爱多彩”唱游世界音乐之旅 纪念专场'
$wsPerform.Range("D19").Value = '中山纪念堂 中山纪念堂'
$wsPerform.Range("E19").Value = '2024.12.20 20:00-12.20 22:00'
$wsPerform.Range("F19").Value = 28
$wsPerform.Range("G19").Value = 380
$wsPerform.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=87739'
$wsPerform.Range("I19").Value = '//i0.hdslb.com/bfs/openplatform/202406/HCPstM8c1718868579079.jpeg'
$wsPerform.Range("B20").Value = '2024-12-24'
$wsPerform.Range("C20").Value = '广州·德国美因茨名家管弦乐团 2025 新年音乐会'
$wsPerform.Range("D20").Value = '人民北路875号（广州市少年宫内） 广州蓓蕾剧院'
$wsPerform.Range("E20").Value = '2024.12.24 19:30-12.24 21:00'
$wsPerform.Range("F20").Value = 0
$wsPerform.Range("G20").Value = 126
$wsPerform.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=93359'
$wsPerform.Range("I20").Value = '//i0.hdslb.com/bfs/openplatform/202410/HaoFdo471728632672864.jpeg'
$wsPerform.Range("B21").Value = '2024-12-29'
$wsPerform.Range("C21").Value = '广州·维也纳皇家交响乐团2025新年音乐会'
$wsPerform.Range("D21").Value = '人民北路696号 广州友谊剧院'
$wsPerform.Range("E21").Value = '2024.12.29 20:00-12.30 21:45'
$wsPerform.Range("F21").Value = 49
$wsPerform.Range("G21").Value = 280
$wsPerform.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=89837'
$wsPerform.Range("I21").Value = '//i2.hdslb.com/bfs/openplatform/202407/OzlirVhz1721882951190.jpeg'
$wsPerform.Range("B22").Value = '2025-01-11'
$wsPerform.Range("C22").Value = '广州·吉冈毅志&高野八诚 2025年见面会'
$wsPerform.Range("D22").Value = '金花街道中山七路333号1906科技圆区3号楼109-1铺、110-1铺、111-1铺 音乐唐人馆'
$wsPerform.Range("E22").Value = '2025.01.11 12:00-01.11 19:00'
$wsPerform.Range("F22").Value = 56
$wsPerform.Range("G22").Value = 480
$wsPerform.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=93488'
$wsPerform.Range("I22").Value = '//i1.hdslb.com/bfs/openplatform/202410/TFuXhHIk1728979778227.jpeg'

# The old last row (23) is now a duplicate of row 22s data; remove it so
# the sheets used range shrinks back from I23 to I22.
$wsPerform.Rows("23:23").Delete()

